# Apply updates to "controle de processos" workbook (database/controle_processos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F (uasg) was stored as text in the whole data range; it is now numeric ---
# Set numeric values for every existing data row (2-36) as well as the new row 37.
$uasgValues = @{
    2  = 787000
    3  = 787310
    4  = 787000
    5  = 787200
    6  = 787900
    7  = 787010
    8  = 787000
    9  = 787000
    10 = 787400
    11 = 787000
    12 = 787010
    13 = 787200
    14 = 787310
    15 = 787000
    16 = 787320
    17 = 787000
    18 = 787000
    19 = 787000
    20 = 787000
    21 = 787000
    22 = 787000
    23 = 787000
    24 = 787400
    25 = 787000
    26 = 787010
    27 = 787310
    28 = 787000
    29 = 787000
    30 = 787200
    31 = 787010
    32 = 787000
    33 = 787000
    34 = 787900
    35 = 787900
    36 = 787000
    37 = 787000
}
foreach ($row in $uasgValues.Keys) {
    $ws.Range("F$row").Value = $uasgValues[$row]
}

# --- num_pregao for row 14 was stored as text "9"; it is now numeric 9 ---
$ws.Range("B14").Value = 9

# --- Row 2: Gêneros Secos process moved from CeIMBra to Com7ºDN ---
$ws.Range("G2").Value = "COMANDO DO 7º DISTRITO NAVAL"
$ws.Range("H2").Value = "Com7ºDN"

# --- Row 5: Material de Pintura moved to Homologado stage ---
$ws.Range("K5").Value = "Homologado"

# --- Row 6: CC 1/2024 object renamed and reassigned to the new CIAB unit ---
$ws.Range("E6").Value = "Nova Sede da CFB"
$ws.Range("G6").Value = "CENTRO DE INSTRUÇÃO E ADESTRAMENTO DE BRASÍLIA"
$ws.Range("H6").Value = "CIAB"
$ws.Range("K6").Value = "Em recurso"

# --- Stage ("etapa") updates for several rows ---
$ws.Range("K7").Value = "Setor Responsável"
$ws.Range("K8").Value = "Edital"
$ws.Range("K9").Value = "Sessão Pública"
$ws.Range("K11").Value = "Sessão Pública"
$ws.Range("K13").Value = "Recomendações AGU"
$ws.Range("K26").Value = "Planejamento"
$ws.Range("K35").Value = "Planejamento"

# --- Row 29: Dispensa4 reassigned from GptFNB to Com7ºDN ---
$ws.Range("G29").Value = "COMANDO DO 7º DISTRITO NAVAL"
$ws.Range("H29").Value = "Com7ºDN"

# --- Rows 34 and 35 (TJIL / TJDL): object and unit updated to the new CIAB entry ---
$ws.Range("E34").Value = "Nova Sede da CFB"
$ws.Range("G34").Value = "CENTRO DE INSTRUÇÃO E ADESTRAMENTO DE BRASÍLIA"
$ws.Range("H34").Value = "CIAB"

$ws.Range("E35").Value = "Nova Sede da CFB"
$ws.Range("G35").Value = "CENTRO DE INSTRUÇÃO E ADESTRAMENTO DE BRASÍLIA"
$ws.Range("H35").Value = "CIAB"

# --- New row 37: second process added under the same 62055.XXXXXX/2024-XX NUP ---
$ws.Range("A37").Value = "PE"
$ws.Range("B37").Value = 22
$ws.Range("C37").Value = 2024
$ws.Range("D37").Value = "62055.XXXXXX/2024-XX"
$ws.Range("G37").Value = "COMANDO DO 7º DISTRITO NAVAL"
$ws.Range("H37").Value = "Com7ºDN"
$ws.Range("K37").Value = "Planejamento"

# --- Update sheet dimension to include the new row ---
$ws.Range("A1:L37").Select | Out-Null
